# Update "negativo" (F) and "positivo" (G) comment counts for several
# news-story rows in the CRONACA comments dataset.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 4;  F = 74; G = 26 },
    @{ Row = 5;  F = 68; G = 32 },
    @{ Row = 6;  F = 66; G = 34 },
    @{ Row = 7;  F = 91; G = 9  },
    @{ Row = 8;  F = 81 },
    @{ Row = 9;  F = 85; G = 15 },
    @{ Row = 10; F = 84; G = 16 },
    @{ Row = 11; F = 84 },
    @{ Row = 15; F = 35; G = 65 },
    @{ Row = 21; F = 88; G = 12 },
    @{ Row = 36; F = 85 }
)

foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, 6).Value = $u.F
    if ($u.ContainsKey('G')) {
        $ws.Cells.Item($u.Row, 7).Value = $u.G
    }
}
